$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030533058377069
$ws.Range("D2").Value = 1.033177462681259
$ws.Range("E2").Value = 1.039552098467655
$ws.Range("F2").Value = 1.048896942223943
$ws.Range("I2").Value = 1.03499182437724
$ws.Range("J2").Value = 1.035673778774571
$ws.Range("K2").Value = 1.035980562479287
$ws.Range("L2").Value = 1.04233696238428
$ws.Range("M2").Value = 1.051655507976379
$ws.Range("N2").Value = 1.015890573910723

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031439371046339
$ws.Range("D3").Value = 1.033997276968464
$ws.Range("E3").Value = 1.040365164007232
$ws.Range("F3").Value = 1.049814137186467
$ws.Range("I3").Value = 1.035137131015873
$ws.Range("J3").Value = 1.036221805680694
$ws.Range("K3").Value = 1.036609273748202
$ws.Range("L3").Value = 1.04296026241975
$ws.Range("M3").Value = 1.052384568828572
$ws.Range("N3").Value = 1.016072226417863

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032025990599636
$ws.Range("D4").Value = 1.034528191861656
$ws.Range("E4").Value = 1.040891795683374
$ws.Range("F4").Value = 1.050408131394257
$ws.Range("I4").Value = 1.035229228095921
$ws.Range("J4").Value = 1.036575967501031
$ws.Range("K4").Value = 1.037015884318485
$ws.Range("L4").Value = 1.043363452152952
$ws.Range("M4").Value = 1.052856197274178
$ws.Range("N4").Value = 1.01618959623993

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032272646116427
$ws.Range("D5").Value = 1.034751492394125
$ws.Range("E5").Value = 1.041113315578746
$ws.Range("F5").Value = 1.050657966569928
$ws.Range("I5").Value = 1.035267483822197
$ws.Range("J5").Value = 1.036724748730428
$ws.Range("K5").Value = 1.037186772472768
$ws.Range("L5").Value = 1.043532921497644
$ws.Range("M5").Value = 1.053054439135875
$ws.Range("N5").Value = 1.016238897064019

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032314063013525
$ws.Range("D6").Value = 1.034788991586366
$ws.Range("E6").Value = 1.041150516968472
$ws.Range("F6").Value = 1.050699921987806
$ws.Range("I6").Value = 1.035273880018659
$ws.Range("J6").Value = 1.036749723365059
$ws.Range("K6").Value = 1.037215462323875
$ws.Range("L6").Value = 1.043561374257221
$ws.Range("M6").Value = 1.053087722974472
$ws.Range("N6").Value = 1.016247172447281

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032029286264422
$ws.Range("D7").Value = 1.034531175207868
$ws.Range("E7").Value = 1.040894755157938
$ws.Range("F7").Value = 1.050411469234487
$ws.Range("I7").Value = 1.035229741086705
$ws.Range("J7").Value = 1.036577955950416
$ws.Range("K7").Value = 1.037018167936018
$ws.Range("L7").Value = 1.043365716736384
$ws.Range("M7").Value = 1.052858846314138
$ws.Range("N7").Value = 1.016190255163317

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030839314129932
$ws.Range("D8").Value = 1.033454430716942
$ws.Range("E8").Value = 1.039826768404568
$ws.Range("F8").Value = 1.049206806500044
$ws.Range("I8").Value = 1.035041329760232
$ws.Range("J8").Value = 1.035859079188069
$ws.Range("K8").Value = 1.036193080665029
$ws.Range("L8").Value = 1.042547635140754
$ws.Range("M8").Value = 1.051901922085161
$ws.Range("N8").Value = 1.015951999500154

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028743819555473
$ws.Range("D9").Value = 1.031560506066065
$ws.Range("E9").Value = 1.037948917060328
$ws.Range("F9").Value = 1.047087993282998
$ws.Range("I9").Value = 1.034694608406145
$ws.Range("J9").Value = 1.034588941468503
$ws.Range("K9").Value = 1.034737634382217
$ws.Range("L9").Value = 1.041105150299586
$ws.Range("M9").Value = 1.050214819844749
$ws.Range("N9").Value = 1.015530866348143

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027347823444923
$ws.Range("D10").Value = 1.030300286966599
$ws.Range("E10").Value = 1.036699844863769
$ws.Range("F10").Value = 1.045678201370149
$ws.Range("I10").Value = 1.034453612145942
$ws.Range("J10").Value = 1.033739970713905
$ws.Range("K10").Value = 1.0337663743173
$ws.Range("L10").Value = 1.040142944837058
$ws.Range("M10").Value = 1.049089573552109
$ws.Range("N10").Value = 1.015249262316032

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026743592553172
$ws.Range("D11").Value = 1.029755184755504
$ws.Range("E11").Value = 1.036159671797509
$ws.Range("F11").Value = 1.045068416853244
$ws.Range("I11").Value = 1.034346931551301
$ws.Range("J11").Value = 1.033371845387852
$ws.Range("K11").Value = 1.033345594321972
$ws.Range("L11").Value = 1.039726183523834
$ws.Range("M11").Value = 1.04860222465575
$ws.Range("N11").Value = 1.015127128119609

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026519192085759
$ws.Range("D12").Value = 1.029552798223707
$ws.Range("E12").Value = 1.035959131574905
$ws.Range("F12").Value = 1.044842017056614
$ws.Range("I12").Value = 1.034306956643097
$ws.Range("J12").Value = 1.033235030836616
$ws.Range("K12").Value = 1.033189266289808
$ws.Range("L12").Value = 1.03957136291829
$ws.Range("M12").Value = 1.048421186363309
$ws.Range("N12").Value = 1.015081732693887

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026567325021118
$ws.Range("D13").Value = 1.029596206785262
$ws.Range("E13").Value = 1.036002143399561
$ws.Range("F13").Value = 1.044890575972989
$ws.Range("I13").Value = 1.034315547177522
$ws.Range("J13").Value = 1.033264381463915
$ws.Range("K13").Value = 1.033222800595355
$ws.Range("L13").Value = 1.039604573210267
$ws.Range("M13").Value = 1.0484600203594
$ws.Range("N13").Value = 1.015091471488917

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026725042763662
$ws.Range("D14").Value = 1.02973845360403
$ws.Range("E14").Value = 1.036143092948199
$ws.Range("F14").Value = 1.045049700500018
$ws.Range("I14").Value = 1.034343634328004
$ws.Range("J14").Value = 1.033360537807332
$ws.Range("K14").Value = 1.033332672836245
$ws.Range("L14").Value = 1.039713386342463
$ws.Range("M14").Value = 1.04858726027694
$ws.Range("N14").Value = 1.015123376316103

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026822222908566
$ws.Range("D15").Value = 1.029826108359675
$ws.Range("E15").Value = 1.036229950454602
$ws.Range("F15").Value = 1.045147755842914
$ws.Range("I15").Value = 1.034360893525722
$ws.Range("J15").Value = 1.033419772808767
$ws.Range("K15").Value = 1.033400364595992
$ws.Range("L15").Value = 1.039780427499164
$ws.Range("M15").Value = 1.048665655012111
$ws.Range("N15").Value = 1.015143030055351

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027387929501834
$ws.Range("D16").Value = 1.030336475966487
$ws.Range("E16").Value = 1.036735708892294
$ws.Range("F16").Value = 1.045718684889322
$ws.Range("I16").Value = 1.034460643206346
$ws.Range("J16").Value = 1.033764391207674
$ws.Range("K16").Value = 1.033794295604828
$ws.Range("L16").Value = 1.040170601481249
$ws.Range("M16").Value = 1.049121915089373
$ws.Range("N16").Value = 1.01525736382468

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027742848529539
$ws.Range("D17").Value = 1.03065677262671
$ws.Range("E17").Value = 1.037053141782681
$ws.Range("F17").Value = 1.046076992585143
$ws.Range("I17").Value = 1.034522590880178
$ws.Range("J17").Value = 1.033980423960027
$ws.Range("K17").Value = 1.034041340567994
$ws.Range("L17").Value = 1.040415315984741
$ws.Range("M17").Value = 1.049408086488185
$ws.Range("N17").Value = 1.015329029683874

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027949890328177
$ws.Range("D18").Value = 1.030843652210705
$ws.Range("E18").Value = 1.03723836086353
$ws.Range("F18").Value = 1.046286051436668
$ws.Range("I18").Value = 1.034558499280747
$ws.Range("J18").Value = 1.034106382414192
$ws.Range("K18").Value = 1.034185416607315
$ws.Range("L18").Value = 1.0405580421607
$ws.Range("M18").Value = 1.049574994740818
$ws.Range("N18").Value = 1.015370812058325

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028020490191413
$ws.Range("D19").Value = 1.030907382768284
$ws.Range("E19").Value = 1.037301526908825
$ws.Range("F19").Value = 1.046357345919299
$ws.Range("I19").Value = 1.034570704994534
$ws.Range("J19").Value = 1.034149322514861
$ws.Range("K19").Value = 1.034234539198363
$ws.Range("L19").Value = 1.040606706087188
$ws.Range("M19").Value = 1.049631904290478
$ws.Range("N19").Value = 1.015385055515841

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027704766634964
$ws.Range("D20").Value = 1.030622402027748
$ws.Range("E20").Value = 1.037019077416551
$ws.Range("F20").Value = 1.04603854292636
$ws.Range("I20").Value = 1.034515967710339
$ws.Range("J20").Value = 1.033957250833392
$ws.Range("K20").Value = 1.03401483715647
$ws.Range("L20").Value = 1.040389061631661
$ws.Range("M20").Value = 1.049377384101008
$ws.Range("N20").Value = 1.01532134258925

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02667859776802
$ws.Range("D21").Value = 1.029696563001414
$ws.Range("E21").Value = 1.036101583933947
$ws.Range("F21").Value = 1.045002839484133
$ws.Range("I21").Value = 1.034335372990238
$ws.Range("J21").Value = 1.033332224262095
$ws.Range("K21").Value = 1.033300319063144
$ws.Range("L21").Value = 1.039681344044519
$ws.Range("M21").Value = 1.048549791702239
$ws.Range("N21").Value = 1.015113981946873

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026033623820921
$ws.Range("D22").Value = 1.029114965147368
$ws.Range("E22").Value = 1.035525322390714
$ws.Range("F22").Value = 1.044352238570764
$ws.Range("I22").Value = 1.034219807149826
$ws.Range("J22").Value = 1.032938803639867
$ws.Range("K22").Value = 1.032850890874417
$ws.Range("L22").Value = 1.039236276701213
$ws.Range("M22").Value = 1.048029364129342
$ws.Range("N22").Value = 1.014983436468839

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026375515595536
$ws.Range("D23").Value = 1.029423231992111
$ws.Range("E23").Value = 1.035830751903177
$ws.Range("F23").Value = 1.044697078301803
$ws.Range("I23").Value = 1.034281261907645
$ws.Range("J23").Value = 1.033147404876582
$ws.Range("K23").Value = 1.033089158160907
$ws.Range("L23").Value = 1.039472224158466
$ws.Range("M23").Value = 1.048305260487125
$ws.Range("N23").Value = 1.015052657051247

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027721974131248
$ws.Range("D24").Value = 1.030637932449179
$ws.Range("E24").Value = 1.037034469434006
$ws.Range("F24").Value = 1.046055916474391
$ws.Range("I24").Value = 1.034518961129822
$ws.Range("J24").Value = 1.033967721925786
$ws.Range("K24").Value = 1.034026812972631
$ws.Range("L24").Value = 1.04040092487982
$ws.Range("M24").Value = 1.049391257221322
$ws.Range("N24").Value = 1.015324816115499

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029285383709943
$ws.Range("D25").Value = 1.032049714764619
$ws.Range("E25").Value = 1.038433894419369
$ws.Range("F25").Value = 1.047635279120597
$ws.Range("I25").Value = 1.034785983423769
$ws.Range("J25").Value = 1.034917697082431
$ws.Range("K25").Value = 1.035114076960747
$ws.Range("L25").Value = 1.041478169063484
$ws.Range("M25").Value = 1.050651072188752
$ws.Range("N25").Value = 1.015639890604769
